$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'310.08"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'0.76%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'37.29"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'-2.12%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'1.71%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'0.07771"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'-1.55%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('B6').Value = 'KuCoinToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D6').Value = "'8.194"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'-0.22%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = "'1.875"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'-7.05%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = "'2.781"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'-11.00%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = "'0.9179"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'-0.90%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.1189"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'-7.15%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1896"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'0.86%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = "'0.09411"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'7.65%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03408"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'-0.52%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09691"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'-0.17%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001373"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'-1.51%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.005819"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'-4.05%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.536"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'-1.06%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = "'4.409"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'1.03%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.3402"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'-1.05%"
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'5.257"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'4.78%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'-1.30%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'2.89%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.02104"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'180.32%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.04335"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'-0.02%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'-1.50%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.004253"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'-7.75%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0001300"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'-63.82%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('D39').Value = "'0.02073"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'-8.50%"
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.05023"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'-0.10%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.007653"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'0.67%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.009830"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'-0.57%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.1345"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'-1.22%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.002170"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'7.01%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.008731"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'2.05%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00006711"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'4.25%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'-0.25%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = 'BOLO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D48').Value = "'0.002934"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'-2.42%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = 'CoinbaseStockToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D49').Value = "'0.001201"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'-0.07%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.00002103"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'-0.25%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.0002002"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'-0.25%"
$ws.Range('E51').Style = 'Normal'
